# Updated variation and parameter storage
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lung")

# Column C is now wider to fit the updated p-value text
$ws.Columns.Item(3).ColumnWidth = 18.0833333333333

# Header cell: day/sample-size label
$ws.Range("A1").Value = "Day 1, n = 1000"

# Oral Dose / AUC_24 row
$ws.Range("C2").Value = 179.07
$ws.Range("D2").Value = 50.1

# Oral Dose / C_max row
$ws.Range("C3").Value = 14.83
$ws.Range("D3").Value = 1.92

# Lung Dose / AUC_24 row
$ws.Range("C4").Value = 184.73
$ws.Range("D4").Value = 48.86

# Lung Dose / C_max row
$ws.Range("C5").Value = 30.3
$ws.Range("D5").Value = 6.26

# Better Dose / AUC_24 row (p-value + effect size)
$ws.Range("C7").Value = "Lung, p = 7.3688e-247"
$ws.Range("D7").Value = 5.66

# Better Dose / C_max row (effect size only; p-value unchanged)
$ws.Range("D8").Value = 15.47
